$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Uniform B:Q values applied to every data row (2-26)
$values = @([double]"0.9999805300082542", [double]"0.9991182316315311", [double]"0.9999999999997352", [double]"0.9999905011343228", [double]"0.9999931658025812", [double]"1.817438266093276e-05", [double]"0.0008230920668076507", [double]"8.756010848999381e-14", [double]"1.357546553520946e-05", [double]"6.787732811384783e-06", [double]"0.0002731515751168351", [double]"0.004263142345844525", [double]"1.000035944600146", [double]"0.004444633141428372", [double]"95.83099500301532", [double]"140.9294005231387")

# New names for column A, keyed by row number
$names = @{
    2 = "model_20_4_0"
    3 = "model_20_4_22"
    4 = "model_20_4_21"
    5 = "model_20_4_20"
    6 = "model_20_4_19"
    7 = "model_20_4_18"
    8 = "model_20_4_17"
    9 = "model_20_4_16"
    10 = "model_20_4_15"
    11 = "model_20_4_14"
    12 = "model_20_4_13"
    13 = "model_20_4_23"
    14 = "model_20_4_12"
    15 = "model_20_4_10"
    16 = "model_20_4_9"
    17 = "model_20_4_8"
    18 = "model_20_4_7"
    19 = "model_20_4_6"
    20 = "model_20_4_5"
    21 = "model_20_4_4"
    22 = "model_20_4_3"
    23 = "model_20_4_2"
    24 = "model_20_4_1"
    25 = "model_20_4_11"
    26 = "model_20_4_24"
}

foreach ($r in $names.Keys) {
    $ws.Cells.Item($r, 1).Value = $names[$r]
    for ($i = 0; $i -lt $values.Count; $i++) {
        $ws.Cells.Item($r, $i + 2).Value = $values[$i]
    }
}
